# The deck's 5th slide ("Exploration:/Obtaining Datasets:/Cleaning:" methodology
# slide) was reordered so it now appears right after slide 4, pushing the four
# research-question slides ("Can the number of restaurants...", "Does a county
# with a higher economic standing...", "Is there a correlation between the
# number of restaurants and obesity...", "Is there a correlation between
# population and number of restaurants...") down by one position each.
#
# i.e. the slide currently at index 9 moves to index 5; no shapes/text/slide
# content are otherwise changed.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$s.MoveTo(5)
